# Generate Report for Archive
#
# Swap the two tracked files' rows (52bcbd4b-....md now sorts after
# 0dd5d789-....md) across the Overview / zh-cn / de-de sheets, and update
# the Status column from "Ready for handoff" to "In Translation" for both
# swapped rows.
#
# Note: the hyperlink relationship targets (rId -> external URL) stay
# attached to the exact same cell refs they always were (rId2 stays on
# A2's GitHub "52bcbd4b...md" URL, rId3 stays on A3's "0dd5d789...md"
# URL, etc.) - only the hyperlink's displayed text is swapped to match
# the new row order. The link target URLs themselves are not touched.

$wb = $excel.ActiveWorkbook

# ---- Overview sheet -------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")

# Original per-ref hyperlink targets (unchanged by this edit).
$urlA2 = "https://github.com/OpenLocalizationTest/oltest/blob/6e43884ed780266b9c5ff4ebe0a6dc449683f2c0/e2e/52bcbd4b-d994-4290-9bba-13a260905a83.md"
$urlA3 = "https://github.com/OpenLocalizationTest/oltest/blob/6e43884ed780266b9c5ff4ebe0a6dc449683f2c0/e2e/0dd5d789-0d2b-468c-9fd2-0ccaf558259d.md"
$urlA4 = "https://github.com/OpenLocalizationTest/oltest/blob/6e43884ed780266b9c5ff4ebe0a6dc449683f2c0/.localization-config"

# Wipe every hyperlink on the sheet (Hyperlinks.Delete on a Range-scoped
# collection clears the whole sheet's collection in this host) so we can
# rebuild the three links cleanly with their new display text.
$ws.Range("A2").Hyperlinks.Delete()

$ws.Hyperlinks.Add($ws.Range("A2"), $urlA2, [Type]::Missing, [Type]::Missing, "0dd5d789-0d2b-468c-9fd2-0ccaf558259d.md")
$ws.Hyperlinks.Add($ws.Range("A3"), $urlA3, [Type]::Missing, [Type]::Missing, "52bcbd4b-d994-4290-9bba-13a260905a83.md")
$ws.Hyperlinks.Add($ws.Range("A4"), $urlA4, [Type]::Missing, [Type]::Missing, ".localization-config")

$ws.Range("B2").Value = "In Translation"
$ws.Range("C2").Value = "In Translation"
$ws.Range("B3").Value = "In Translation"
$ws.Range("C3").Value = "In Translation"

# ---- zh-cn sheet ------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")

$urlA2 = "https://github.com/OpenLocalizationTest/oltest/blob/6e43884ed780266b9c5ff4ebe0a6dc449683f2c0/e2e/52bcbd4b-d994-4290-9bba-13a260905a83.md"
$urlC2 = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/243e4f2b3f1d1ba88ad213aee5513dbe58c54a03/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/52bcbd4b-d994-4290-9bba-13a260905a83.5a3aff9d203e6d6cb7e2f7feb8bc6241b16153d9.zh-cn.xlf"
$urlA3 = "https://github.com/OpenLocalizationTest/oltest/blob/6e43884ed780266b9c5ff4ebe0a6dc449683f2c0/e2e/0dd5d789-0d2b-468c-9fd2-0ccaf558259d.md"
$urlC3 = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/243e4f2b3f1d1ba88ad213aee5513dbe58c54a03/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/0dd5d789-0d2b-468c-9fd2-0ccaf558259d.7cfd6499ef4b5ef1dbd128f298bb1e949efd6c25.zh-cn.xlf"
$urlA4 = "https://github.com/OpenLocalizationTest/oltest/blob/6e43884ed780266b9c5ff4ebe0a6dc449683f2c0/.localization-config"

$ws.Range("A2").Hyperlinks.Delete()

$ws.Hyperlinks.Add($ws.Range("A2"), $urlA2, [Type]::Missing, [Type]::Missing, "0dd5d789-0d2b-468c-9fd2-0ccaf558259d.md")
$ws.Hyperlinks.Add($ws.Range("C2"), $urlC2, [Type]::Missing, [Type]::Missing, "0dd5d789-0d2b-468c-9fd2-0ccaf558259d.7cfd6499ef4b5ef1dbd128f298bb1e949efd6c25.zh-cn.xlf")
$ws.Hyperlinks.Add($ws.Range("A3"), $urlA3, [Type]::Missing, [Type]::Missing, "52bcbd4b-d994-4290-9bba-13a260905a83.md")
$ws.Hyperlinks.Add($ws.Range("C3"), $urlC3, [Type]::Missing, [Type]::Missing, "52bcbd4b-d994-4290-9bba-13a260905a83.5a3aff9d203e6d6cb7e2f7feb8bc6241b16153d9.zh-cn.xlf")
$ws.Hyperlinks.Add($ws.Range("A4"), $urlA4, [Type]::Missing, [Type]::Missing, ".localization-config")

$ws.Range("B2").Value = "In Translation"
$ws.Range("B3").Value = "In Translation"

# ---- de-de sheet --------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")

$urlA2 = "https://github.com/OpenLocalizationTest/oltest/blob/6e43884ed780266b9c5ff4ebe0a6dc449683f2c0/e2e/52bcbd4b-d994-4290-9bba-13a260905a83.md"
$urlC2 = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/ac8e341a0ca65744383c927fbaebfbf02bf6b4a9/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/52bcbd4b-d994-4290-9bba-13a260905a83.5a3aff9d203e6d6cb7e2f7feb8bc6241b16153d9.de-de.xlf"
$urlA3 = "https://github.com/OpenLocalizationTest/oltest/blob/6e43884ed780266b9c5ff4ebe0a6dc449683f2c0/e2e/0dd5d789-0d2b-468c-9fd2-0ccaf558259d.md"
$urlC3 = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/ac8e341a0ca65744383c927fbaebfbf02bf6b4a9/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/0dd5d789-0d2b-468c-9fd2-0ccaf558259d.7cfd6499ef4b5ef1dbd128f298bb1e949efd6c25.de-de.xlf"
$urlA4 = "https://github.com/OpenLocalizationTest/oltest/blob/6e43884ed780266b9c5ff4ebe0a6dc449683f2c0/.localization-config"

$ws.Range("A2").Hyperlinks.Delete()

$ws.Hyperlinks.Add($ws.Range("A2"), $urlA2, [Type]::Missing, [Type]::Missing, "0dd5d789-0d2b-468c-9fd2-0ccaf558259d.md")
$ws.Hyperlinks.Add($ws.Range("C2"), $urlC2, [Type]::Missing, [Type]::Missing, "0dd5d789-0d2b-468c-9fd2-0ccaf558259d.7cfd6499ef4b5ef1dbd128f298bb1e949efd6c25.de-de.xlf")
$ws.Hyperlinks.Add($ws.Range("A3"), $urlA3, [Type]::Missing, [Type]::Missing, "52bcbd4b-d994-4290-9bba-13a260905a83.md")
$ws.Hyperlinks.Add($ws.Range("C3"), $urlC3, [Type]::Missing, [Type]::Missing, "52bcbd4b-d994-4290-9bba-13a260905a83.5a3aff9d203e6d6cb7e2f7feb8bc6241b16153d9.de-de.xlf")
$ws.Hyperlinks.Add($ws.Range("A4"), $urlA4, [Type]::Missing, [Type]::Missing, ".localization-config")

$ws.Range("B2").Value = "In Translation"
$ws.Range("B3").Value = "In Translation"
